$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6276753333333334
$ws.Range("H2").Value = 1.883026
$ws.Range("I2").Value = 0.01507055680360577
$ws.Range("J2").Value = 0.01507055680360577
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1612466666666667
$ws.Range("N2").Value = 0.48374
$ws.Range("O2").Value = 0.04266180225345401
$ws.Range("P2").Value = 0.04266180225345401
$ws.Range("Q2").Value = 0.1012105552488889
$ws.Range("R2").Value = 0.91089499724
$ws.Range("S2").Value = 0.0006429371142048753
$ws.Range("T2").Value = 0.0006429371142048753
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6276753333333334
$ws.Range("H3").Value = 1.883026
$ws.Range("I3").Value = 0.01507055680360577
$ws.Range("J3").Value = 0.01507055680360577
$ws.Range("O3").Value = 0.1964844360795697
$ws.Range("P3").Value = 0.1964844360795697
$ws.Range("Q3").Value = 0.4661382741224445
$ws.Range("R3").Value = 4.195244467102
$ws.Range("S3").Value = 0.002961129854961601
$ws.Range("T3").Value = 0.002961129854961601
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6276753333333334
$ws.Range("H4").Value = 1.883026
$ws.Range("I4").Value = 0.01507055680360577
$ws.Range("J4").Value = 0.01507055680360577
$ws.Range("O4").Value = 0.7608537616669764
$ws.Range("P4").Value = 0.7608537616669764
$ws.Range("Q4").Value = 1.805044035036889
$ws.Range("R4").Value = 16.245396315332
$ws.Range("S4").Value = 0.01146648983443929
$ws.Range("T4").Value = 0.01146648983443929
$ws.Range("I5").Value = 0.6396010460118555
$ws.Range("J5").Value = 0.6396010460118555
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1612466666666667
$ws.Range("N5").Value = 0.48374
$ws.Range("O5").Value = 0.04266180225345401
$ws.Range("P5").Value = 0.04266180225345401
$ws.Range("Q5").Value = 4.295420391444445
$ws.Range("R5").Value = 38.658783523
$ws.Range("S5").Value = 0.02728653334606012
$ws.Range("T5").Value = 0.02728653334606012
$ws.Range("I6").Value = 0.6396010460118555
$ws.Range("J6").Value = 0.6396010460118555
$ws.Range("O6").Value = 0.1964844360795697
$ws.Range("P6").Value = 0.1964844360795697
$ws.Range("S6").Value = 0.1256716508415423
$ws.Range("T6").Value = 0.1256716508415423
$ws.Range("I7").Value = 0.6396010460118555
$ws.Range("J7").Value = 0.6396010460118555
$ws.Range("O7").Value = 0.7608537616669764
$ws.Range("P7").Value = 0.7608537616669764
$ws.Range("S7").Value = 0.4866428618242531
$ws.Range("T7").Value = 0.4866428618242531
$ws.Range("I8").Value = 0.3453283971845387
$ws.Range("J8").Value = 0.3453283971845388
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1612466666666667
$ws.Range("N8").Value = 0.48374
$ws.Range("O8").Value = 0.04266180225345401
$ws.Range("P8").Value = 0.04266180225345401
$ws.Range("Q8").Value = 2.319149801677778
$ws.Range("R8").Value = 20.8723482151
$ws.Range("S8").Value = 0.01473233179318901
$ws.Range("T8").Value = 0.01473233179318902
$ws.Range("I9").Value = 0.3453283971845387
$ws.Range("J9").Value = 0.3453283971845388
$ws.Range("O9").Value = 0.1964844360795697
$ws.Range("P9").Value = 0.1964844360795697
$ws.Range("R9").Value = 96.13029342585502
$ws.Range("S9").Value = 0.06785165538306574
$ws.Range("T9").Value = 0.06785165538306576
$ws.Range("I10").Value = 0.3453283971845387
$ws.Range("J10").Value = 0.3453283971845388
$ws.Range("O10").Value = 0.7608537616669764
$ws.Range("P10").Value = 0.7608537616669764
$ws.Range("Q10").Value = 41.36097767254778
$ws.Range("R10").Value = 372.24879905293
$ws.Range("S10").Value = 0.262744410008284
$ws.Range("T10").Value = 0.262744410008284
